# Weekly driver report update for 2025-04-19
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Critical Minutes values for the bad-driver row and totals row
$ws.Range("C3").Value = 175
$ws.Range("C4").Value = 175

# Clear the Driver Vintage date for the first good-driver row (E12)
$ws.Range("E12").ClearContents()

# Update Total Samples for the third good-driver row (B14)
$ws.Range("B14").Value = 265400
